# Hortaliza, Vega Monumental Concepción - Cebolla
# Insert two new weekly price rows (1a/2a cosecha, 23-Dec-2022) above the
# existing row 566 ("1a (guarda)" / 44824), shifting the rest of the table
# down by two rows (566-639 -> 568-641).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 566, pushing rows 566:639 down to 568:641.
# -4121 = xlShiftDown
$ws.Range("A566:R567").EntireRow.Insert(-4121)

# Row 566: 1a (cosecha)
$ws.Range("A566").Value = 11
$ws.Range("B566").Value = "Vega Monumental Concepción"
$ws.Range("C566").Value = "Bíobío"
$ws.Range("D566").Value = 44918
$ws.Range("D566").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E566").Value = 8
$ws.Range("F566").Value = 100112004
$ws.Range("G566").Value = "Cebolla"
$ws.Range("H566").Value = "Sin especificar"
$ws.Range("I566").Value = "1a (cosecha)"
$ws.Range("J566").Value = 400
$ws.Range("K566").Value = 12000
$ws.Range("L566").Value = 12000
$ws.Range("M566").Value = 12000
$ws.Range("N566").Value = "`$/malla 18 kilos"
$ws.Range("O566").Value = "Región de O'Higgins"
$ws.Range("P566").Value = 667
$ws.Range("Q566").Value = 18
$ws.Range("R566").Value = "Hortaliza"

# Row 567: 2a (cosecha)
$ws.Range("A567").Value = 11
$ws.Range("B567").Value = "Vega Monumental Concepción"
$ws.Range("C567").Value = "Bíobío"
$ws.Range("D567").Value = 44918
$ws.Range("D567").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E567").Value = 8
$ws.Range("F567").Value = 100112004
$ws.Range("G567").Value = "Cebolla"
$ws.Range("H567").Value = "Sin especificar"
$ws.Range("I567").Value = "2a (cosecha)"
$ws.Range("J567").Value = 400
$ws.Range("K567").Value = 10000
$ws.Range("L567").Value = 10000
$ws.Range("M567").Value = 10000
$ws.Range("N567").Value = "`$/malla 18 kilos"
$ws.Range("O567").Value = "Región de O'Higgins"
$ws.Range("P567").Value = 556
$ws.Range("Q567").Value = 18
$ws.Range("R567").Value = "Hortaliza"
